$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '29.268.74'
Set-TextValue $ws 'E2' '  +0.47%  '
Set-TextValue $ws 'D3' '1.858.44'
Set-TextValue $ws 'E3' '  +0.35%  '
Set-TextValue $ws 'D4' '1.001'
Set-TextValue $ws 'E5' '  +2.33%  '
Set-TextValue $ws 'D6' '238.39'
Set-TextValue $ws 'E6' '  +0.24%  '
Set-TextValue $ws 'D7' '1.000'
Set-TextValue $ws 'E7' '  +0.04%  '
Set-TextValue $ws 'D8' '0.08001'
Set-TextValue $ws 'E8' '  +2.77%  '
Set-TextValue $ws 'D9' '0.3030'
Set-TextValue $ws 'E9' '  -0.18%  '
Set-TextValue $ws 'D10' '23.52'
Set-TextValue $ws 'E10' '  +1.30%  '
Set-TextValue $ws 'D11' '0.08193'
Set-TextValue $ws 'E11' '  +0.73%  '
Set-TextValue $ws 'D12' '1.878.10'
Set-TextValue $ws 'E12' '  +1.46%  '
Set-TextValue $ws 'D13' '5.205'
Set-TextValue $ws 'E13' '  +0.08%  '
Set-TextValue $ws 'D14' '0.7077'
Set-TextValue $ws 'E14' '  -2.32%  '
Set-TextValue $ws 'D15' '89.77'
Set-TextValue $ws 'E15' '  +0.84%  '
Set-TextValue $ws 'D16' '29.429.14'
Set-TextValue $ws 'E16' '  +0.98%  '
Set-TextValue $ws 'D17' '5.840'
Set-TextValue $ws 'E17' '  +1.81%  '
Set-TextValue $ws 'D18' '0.000007923'
Set-TextValue $ws 'E18' '  +1.30%  '
Set-TextValue $ws 'D19' '13.30'
Set-TextValue $ws 'E19' '  +0.86%  '
Set-TextValue $ws 'D20' '238.28'
Set-TextValue $ws 'E20' '  +1.10%  '
Set-TextValue $ws 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D21' '2.165.43'
Set-TextValue $ws 'E21' '  +2.34%  '
Set-TextValue $ws 'B22' 'Dai'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D22' '1.002'
Set-TextValue $ws 'E22' '  +0.24%  '
Set-TextValue $ws 'E23' '  +0.06%  '
Set-TextValue $ws 'D24' '7.495'
Set-TextValue $ws 'E24' '  -1.33%  '
Set-TextValue $ws 'D25' '163.08'
Set-TextValue $ws 'E25' '  +1.01%  '
Set-TextValue $ws 'D26' '8.894'
Set-TextValue $ws 'E26' '  -0.71%  '
Set-TextValue $ws 'D27' '0.1441'
Set-TextValue $ws 'E27' '  +1.49%  '
Set-TextValue $ws 'D28' '18.11'
Set-TextValue $ws 'E28' '  +0.30%  '
Set-TextValue $ws 'D29' '1.927'
Set-TextValue $ws 'E29' '  -1.99%  '
Set-TextValue $ws 'D30' '1.426'
Set-TextValue $ws 'E30' '  +2.07%  '
Set-TextValue $ws 'D31' '1.476'
Set-TextValue $ws 'E31' '  -0.61%  '
Set-TextValue $ws 'D32' '4.379'
Set-TextValue $ws 'E32' '  -3.19%  '
Set-TextValue $ws 'D33' '4.031'
Set-TextValue $ws 'E33' '  +0.85%  '
Set-TextValue $ws 'D34' '0.05198'
Set-TextValue $ws 'E34' '  +0.04%  '
Set-TextValue $ws 'D35' '1.163'
Set-TextValue $ws 'E35' '  -1.33%  '
Set-TextValue $ws 'D36' '0.7185'
Set-TextValue $ws 'E36' '  +2.19%  '
Set-TextValue $ws 'D37' '1.003'
Set-TextValue $ws 'E37' '  -2.32%  '
Set-TextValue $ws 'D38' '2.676'
Set-TextValue $ws 'E38' '  +0.88%  '
Set-TextValue $ws 'D39' '0.01854'
Set-TextValue $ws 'E39' '  +0.09%  '
Set-TextValue $ws 'D40' '2.724'
Set-TextValue $ws 'E40' '  +1.91%  '
Set-TextValue $ws 'D41' '0.9395'
Set-TextValue $ws 'E41' '  +3.72%  '
Set-TextValue $ws 'D42' '1.158.94'
Set-TextValue $ws 'E42' '  +5.33%  '
Set-TextValue $ws 'D43' '6.013'
Set-TextValue $ws 'E43' '  -0.14%  '
Set-TextValue $ws 'D44' '0.4273'
Set-TextValue $ws 'E44' '  +0.11%  '
Set-TextValue $ws 'D45' '70.77'
Set-TextValue $ws 'E45' '  +0.45%  '
Set-TextValue $ws 'D46' '1.000'
Set-TextValue $ws 'E46' '  +0.04%  '
Set-TextValue $ws 'D47' '102.98'
Set-TextValue $ws 'E47' '  -0.17%  '
Set-TextValue $ws 'D48' '0.5297'
Set-TextValue $ws 'E48' '  -4.33%  '
Set-TextValue $ws 'B49' 'RocketPoolETH'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws 'D49' '2.043.88'
Set-TextValue $ws 'E49' '  +2.07%  '
Set-TextValue $ws 'B50' 'RenderToken'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D50' '1.766'
Set-TextValue $ws 'E50' '  +0.47%  '
Set-TextValue $ws 'D51' '9.170'
Set-TextValue $ws 'E51' '  +0.31%  '
